$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Recorded By" ordering swap: "System, dnasr281@gmail.com"
#    -> "dnasr281@gmail.com, System"  (many rows across the sheet)
# ---------------------------------------------------------------------------
$swapCells = @("G2","G20","G22","G39","G41","G58","G60","G77","G95","G113", `
               "G131","G149","G167","G169","G186","G188","G205","G207")
foreach ($c in $swapCells) {
    $ws.Range($c).Value = "dnasr281@gmail.com, System"
}

# ---------------------------------------------------------------------------
# 2) Summary statistics block (K/L columns) near the top of the sheet
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 72
$ws.Range("L7").Value = 0
$ws.Range("L9").Value = "32.4%"
$ws.Range("L10").Value = "77.4%"

# ---------------------------------------------------------------------------
# 3) Per-group statistics block (columns O/P/Q/R/S) for rows 16-18 and 24-26
# ---------------------------------------------------------------------------
$ws.Range("O16").Value = 6
$ws.Range("P16").Value = 0
$ws.Range("R16").Value = "31.6%"
$ws.Range("S16").Value = "74.2%"

$ws.Range("O17").Value = 6
$ws.Range("P17").Value = 0
$ws.Range("R17").Value = "31.6%"
$ws.Range("S17").Value = "57.0%"

$ws.Range("O18").Value = 6
$ws.Range("P18").Value = 0
$ws.Range("R18").Value = "31.6%"
$ws.Range("S18").Value = "82.5%"

$ws.Range("O24").Value = 6
$ws.Range("P24").Value = 0
$ws.Range("R24").Value = "31.6%"
$ws.Range("S24").Value = "66.7%"

$ws.Range("O25").Value = 6
$ws.Range("P25").Value = 0
$ws.Range("R25").Value = "31.6%"
$ws.Range("S25").Value = "73.6%"

$ws.Range("O26").Value = 6
$ws.Range("P26").Value = 0
$ws.Range("R26").Value = "31.6%"
$ws.Range("S26").Value = "68.4%"

# ---------------------------------------------------------------------------
# 4) Session rows that flipped from "Not Recorded" (pending/missing, pink
#    highlight) to "Recorded" (green highlight) once attendance was taken.
#    Copying the formatting of an existing "Recorded" row (row 2) reuses the
#    same underlying style instead of minting a brand-new one.
# ---------------------------------------------------------------------------
$recordedFormatSource = $ws.Range("A2:I2")

$flippedRows = @(
    @{ Row = 24;  Group = "B1-10"; Students = "29/31" },
    @{ Row = 43;  Group = "B1-11"; Students = "14/19" },
    @{ Row = 62;  Group = "B1-12"; Students = "17/21" },
    @{ Row = 171; Group = "B1-7";  Students = "27/27" },
    @{ Row = 190; Group = "B1-8";  Students = "28/29" },
    @{ Row = 209; Group = "B1-9";  Students = "29/29" }
)

foreach ($item in $flippedRows) {
    $r = $item.Row
    $recordedFormatSource.Copy($ws.Range("A" + $r + ":I" + $r))

    $ws.Range("A$r").Value = "Year 5"
    $ws.Range("B$r").Value = $item.Group
    $ws.Range("C$r").Value = "GENERAL SURGERY"
    $ws.Range("D$r").Value = "5"
    $ws.Range("E$r").Value = "23/12/2025"
    $ws.Range("F$r").Value = "12:00:00"
    $ws.Range("G$r").Value = "dnasr281@gmail.com, System"
    $ws.Range("H$r").Value = $item.Students
    $ws.Range("I$r").Value = "Recorded"
}

# ---------------------------------------------------------------------------
# 5) Column I width (column 9) shrinks from 14 to 10
# ---------------------------------------------------------------------------
$ws.Columns.Item(9).ColumnWidth = 10
